# Insert 6 new rows before row 585, shifting existing rows 585-594 down to 591-600.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A585:A590").EntireRow.Insert()

# Common column values shared across all data rows in this sheet.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112045
$categoria = "Zapallo"
$variedad  = "Camote"
$unidad    = "`$/kilo (volumen en unidades)"
$kgUnid    = 1
$clasif    = "Hortaliza"

# Data for the newly inserted rows (585-590).
$newRows = @(
    @{ Row=585; Fecha=44448; Calidad="1a (guarda)";    Volumen=160; PMin=600; PMax=650; PProm=625; Origen="Provincia de Maipo" },
    @{ Row=586; Fecha=44448; Calidad="1a (guarda)";    Volumen=106; PMin=600; PMax=650; PProm=625; Origen="Región de O'Higgins" },
    @{ Row=587; Fecha=44448; Calidad="1a nueva(o)";    Volumen=97;  PMin=600; PMax=650; PProm=625; Origen="Perú" },
    @{ Row=588; Fecha=44448; Calidad="2a (guarda)";    Volumen=97;  PMin=450; PMax=500; PProm=475; Origen="Provincia de Maipo" },
    @{ Row=589; Fecha=44448; Calidad="2a (guarda)";    Volumen=79;  PMin=450; PMax=500; PProm=475; Origen="Región de O'Higgins" },
    @{ Row=590; Fecha=44448; Calidad="2a nueva(o)";    Volumen=52;  PMin=450; PMax=500; PProm=475; Origen="Perú" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $kgUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
